# Update Leave Card 5/22/2023 1:34 PM
# Insert two new leave entries (VL(2-0-0) on 5/2,3/2023 and two SL(3-0-0) entries
# on 4/25,27,28/2023 and 5/5,8,9/2023) into the leave card table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# --- 1. Make room: insert two new rows into the table body (at the row that
#        currently holds the 6/1/2023 period), pushing everything below down. ---
$ws.Range("96:97").Insert()

# The insert copies default formatting into the new rows; restore the normal
# table-row formatting by pasting the formats from the row right above.
$ws.Range("A95:K95").Copy()
$ws.Range("A96:K97").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Re-apply the calculated "EARNED " helper-column formula on the two new rows
# (PasteSpecial only copied formats, not formulas).
$ws.Range("G96").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G97").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 2. Grow the table definition to cover the two new rows. ---
$tbl.Resize($ws.Range("A8:K143"))

# Resizing the table re-propagates the calculated column formula to the
# trailing rows using a syntax this engine mis-evaluates; rewrite it using the
# same structured-reference form used everywhere else in the column.
$ws.Range("G142").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G143").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 3. Existing row 94 (period 4/1/2023): record 1.25 days EARNED. ---
$ws.Range("C94").Value = 1.25

# --- 4. Row 95 (period 5/1/2023): VL(2-0-0) for 2 days, 5/2,3/2023. ---
$ws.Range("B95").Value = "VL(2-0-0)"
$ws.Range("D95").Value = 2
$ws.Range("K95").Value = "5/2,3/2023"

# --- 5. New row 96: SL(3-0-0) for 3 days, 4/25,27,28/2023. ---
$ws.Range("B96").Value = "SL(3-0-0)"
$ws.Range("H96").Value = 3
$ws.Range("K96").Value = "4/25,27,28/2023"

# --- 6. New row 97: SL(3-0-0) for 3 days, 5/5,8,9/2023. ---
$ws.Range("B97").Value = "SL(3-0-0)"
$ws.Range("H97").Value = 3
$ws.Range("K97").Value = "5/5,8,9/2023"

# --- 7. Update the remembered selection to match the author's last edit. ---
[void]$ws.Activate()
[void]$ws.Range("H98").Select()
